# Rerun w/ tables and flipped april
# The "4" worksheet (April) holds PC timeseries data in A2:C43.
# Flip the sign of every numeric value in that range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("4")

$rng = $ws.Range("A2:C43")
foreach ($cell in $rng.Cells) {
    $v = $cell.Value2
    if ($v -ne $null) {
        $cell.Value2 = -1 * $v
    }
}
